$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

$ws.Range('D2').Value = '27.734.04'
$ws.Range('E2').Value = '  +1.33%  '
$ws.Range('D3').Value = '1.876.27'
$ws.Range('E3').Value = '  +1.55%  '
Set-TextValue 'D4' '1.004'
$ws.Range('E4').Value = '  +0.08%  '
Set-TextValue 'D5' '332.76'
$ws.Range('E5').Value = '  +3.61%  '
$ws.Range('E6').Value = '  +0.10%  '
Set-TextValue 'D7' '0.4734'
$ws.Range('E7').Value = '  +6.13%  '
Set-TextValue 'D8' '0.3961'
$ws.Range('E8').Value = '  +3.53%  '
Set-TextValue 'D9' '47.70'
$ws.Range('E9').Value = '  -2.75%  '
Set-TextValue 'D10' '0.08044'
$ws.Range('E10').Value = '  +2.92%  '
$ws.Range('E11').Value = '  +1.01%  '
Set-TextValue 'D12' '21.90'
$ws.Range('E12').Value = '  +2.37%  '
$ws.Range('D13').Value = '1.875.10'
$ws.Range('E13').Value = '  +1.19%  '
Set-TextValue 'D14' '5.961'
$ws.Range('E14').Value = '  +2.17%  '
Set-TextValue 'D15' '7.166'
$ws.Range('E15').Value = '  +1.11%  '
Set-TextValue 'D16' '1.006'
$ws.Range('E16').Value = '  +0.08%  '
Set-TextValue 'D17' '0.00001051'
$ws.Range('E17').Value = '  +2.90%  '
Set-TextValue 'D18' '87.27'
$ws.Range('E18').Value = '  +2.51%  '
Set-TextValue 'D19' '0.06635'
$ws.Range('E19').Value = '  +2.12%  '
Set-TextValue 'D20' '17.25'
$ws.Range('E20').Value = '  +1.94%  '
Set-TextValue 'D21' '1.003'
$ws.Range('E21').Value = '  +0.20%  '
$ws.Range('D22').Value = '27.733.26'
$ws.Range('E22').Value = '  +1.32%  '
$ws.Range('E23').Value = '  +0.67%  '
Set-TextValue 'D24' '11.05'
$ws.Range('E24').Value = '  +2.80%  '
Set-TextValue 'D25' '2.298'
$ws.Range('E25').Value = '  +1.68%  '
$ws.Range('D26').Value = '2.100.07'
$ws.Range('E26').Value = '  +1.54%  '
Set-TextValue 'D27' '156.68'
$ws.Range('E27').Value = '  +3.46%  '
$ws.Range('E28').Value = '  +4.77%  '
Set-TextValue 'D29' '2.105'
$ws.Range('E29').Value = '  +2.68%  '
Set-TextValue 'D30' '5.584'
$ws.Range('E30').Value = '  +1.80%  '
Set-TextValue 'D31' '122.60'
$ws.Range('E31').Value = '  +2.26%  '
$ws.Range('E32').Value = '  +4.97%  '
Set-TextValue 'D33' '0.09562'
$ws.Range('E33').Value = '  +2.69%  '
Set-TextValue 'D34' '1.455'
$ws.Range('E34').Value = '  -0.85%  '
Set-TextValue 'D35' '3.634'
$ws.Range('E35').Value = '  +0.16%  '
Set-TextValue 'D36' '5.306'
$ws.Range('E36').Value = '  +1.92%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D37' '0.02267'
$ws.Range('E37').Value = '  +2.41%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D38' '0.06113'
$ws.Range('E38').Value = '  +2.94%  '
Set-TextValue 'D39' '1.228'
$ws.Range('E39').Value = '  +1.88%  '
Set-TextValue 'D40' '8.187'
$ws.Range('E40').Value = '  -1.56%  '
Set-TextValue 'D41' '1.002'
$ws.Range('E41').Value = '  +0.12%  '
Set-TextValue 'D42' '0.5999'
$ws.Range('E42').Value = '  +1.61%  '
Set-TextValue 'D43' '0.1914'
$ws.Range('E43').Value = '  +3.62%  '
Set-TextValue 'D44' '10.27'
$ws.Range('E44').Value = '  +0.18%  '
$ws.Range('B45').Value = 'WEMIXTOKEN'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 'D45' '1.268'
$ws.Range('E45').Value = '  +0.75%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue 'D46' '0.5719'
$ws.Range('E46').Value = '  +1.16%  '
Set-TextValue 'D47' '12.24'
$ws.Range('E47').Value = '  +0.61%  '
Set-TextValue 'D48' '3.413'
$ws.Range('E48').Value = '  +1.68%  '
Set-TextValue 'D49' '1.936'
$ws.Range('E49').Value = '  +1.24%  '
Set-TextValue 'D50' '0.06807'
$ws.Range('E50').Value = '  -0.80%  '
